# Add a new "Spain" market sheet, modelled on the existing "Italy" sheet,
# and populate it with the Spain/Zettler test data.

$wb = $excel.ActiveWorkbook

$italy = $wb.Worksheets.Item("Italy")

# Before duplicating it, normalise the Italy sheet's own selection to the
# full used range (this mirrors what happens in the source workbook once
# the new sheet becomes the active tab).
$italy.Select()
$italy.Range("A1:D14").Select()

# Duplicate "Italy" -> new sheet placed right after it; this carries over
# all formatting, merged cells, styles, page setup, etc.
$italy.Copy([System.Reflection.Missing]::Value, $italy)
$spain = $wb.Worksheets.Item($wb.Worksheets.Count)
$spain.Name = "Spain"

# Update the market name and printer model for Spain.
$spain.Range("B2").Value = "Spain Market"
$spain.Range("B4").Value = "NGC-3103/T2056"

# The new text is a different length than "Italy Market" / "NGC-3145/T2446",
# so the columns were re-fitted to the new content.
$spain.Columns.Item(1).ColumnWidth = 24.3
$spain.Columns.Item(2).ColumnWidth = 14.3
$spain.Columns.Item(4).ColumnWidth = 21.3

# Make the new sheet the active tab, with C12 selected.
$spain.Select()
$spain.Range("C12").Select()
